# Doing Updates for Financials
#
# The FDP sheet lists yearly Income Statement / Balance Sheet / Cash Flow
# figures, one column per fiscal year, newest first starting at column D.
# This update adds a brand-new "latest year" column at D (period ending
# 2018-12-29, serial 43462) and pushes all the previously-reported years
# one column to the right (D->E, E->F, ... K->L).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FDP")

# Insert a new blank column at D; existing D:K data shifts to E:L.
$ws.Columns("D:D").Insert()

# The freshly inserted column has no formatting of its own yet - bring
# over the number formats (date / #,##0) from the column that used to be
# D and now sits at E, for every row that actually carries data (7-102).
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---- Income Statement --------------------------------------------------
$ws.Range("D7").Value = 43462
$ws.Range("D8").Value = 4493900
$ws.Range("D9").Value = 4214100
$ws.Range("D10").Value = 279800
$ws.Range("D12").Value = 3200
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 53600
$ws.Range("D15").Value = 0
$ws.Range("D17").Value = 4455300
$ws.Range("D18").Value = 38600
$ws.Range("D20").Value = -14800
$ws.Range("D21").Value = 124300
$ws.Range("D22").Value = 23600
$ws.Range("D23").Value = 200
$ws.Range("D24").Value = 16100
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = -15900
$ws.Range("D27").Value = -21900
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 14800
$ws.Range("D33").Value = -21900
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = -21900

# ---- Balance Sheet -------------------------------------------------------
$ws.Range("D38").Value = 43462
$ws.Range("D41").Value = 21300
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 473500
$ws.Range("D44").Value = 565300
$ws.Range("D45").Value = 78700
$ws.Range("D46").Value = 1138800
$ws.Range("D47").Value = 6100
$ws.Range("D48").Value = 1392200
$ws.Range("D49").Value = 590300
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 127800
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 3255200
$ws.Range("D57").Value = 330000
$ws.Range("D58").Value = 500
$ws.Range("D59").Value = 255500
$ws.Range("D60").Value = 586000
$ws.Range("D61").Value = 661900
$ws.Range("D62").Value = 237700
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 1563200
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 1206000
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 1692000
$ws.Range("D77").Value = 0

# ---- Cash Flow Statement --------------------------------------------------
$ws.Range("D80").Value = 43462
$ws.Range("D81").Value = -21900
$ws.Range("D83").Value = 100500
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 246600
$ws.Range("D91").Value = -150500
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -494800
$ws.Range("D96").Value = -29000
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 242000
$ws.Range("D101").Value = 2400
$ws.Range("D102").Value = -3800
